$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 86.29678392075563
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 198602092.4331654
